$d = $word.ActiveDocument

# The replacement OOXML fragment for the paragraph's run content: splits
# "Learnyounode solution" into 3 runs, with "Learnyounode" wrapped in
# spellcheck proofErr markers, matching the target diff.
$fragmentTemplate = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Learnyounode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> sol</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>ution</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

# Find every paragraph whose text is exactly "Book solution" (bold heading
# runs) and replace just the run content (not the paragraph mark), so the
# paragraph's own properties (<w:pPr>, rsid attributes, etc.) are untouched.
foreach ($p in $d.Paragraphs) {
    $pText = $p.Range.Text
    if ($pText -eq "Book solution`r") {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        # Exclude the trailing paragraph mark from the replaced range.
        $runRange = $d.Range($pStart, $pEnd - 1)
        $runRange.InsertXML($fragmentTemplate)
    }
}
